$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 88030
$ws.Range("B2").Value = "Srta. Maria Clara Aparecida"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("D2").Value = "Doenca"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45105
$ws.Range("G2").Value = 4549.04

$ws.Range("A3").Value = 52654
$ws.Range("B3").Value = "Dra. Sophia Ramos"
$ws.Range("C3").Value = "Operacoes"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45093
$ws.Range("G3").Value = 8637.85

$ws.Range("A4").Value = 31885
$ws.Range("B4").Value = "Bruna Duarte"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45087
$ws.Range("G4").Value = 2897.22

$ws.Range("A5").Value = 51460
$ws.Range("B5").Value = "Maria Camargo"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45099
$ws.Range("G5").Value = 5883.13

$ws.Range("A6").Value = 91193
$ws.Range("B6").Value = "Sr. Davi Miguel da Conceição"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45086
$ws.Range("G6").Value = 4715.91

$ws.Range("A7").Value = 60470
$ws.Range("B7").Value = "Sr. Apollo da Rosa"
$ws.Range("C7").Value = "P&D"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45103
$ws.Range("G7").Value = 8345.23

$ws.Range("A8").Value = 52489
$ws.Range("B8").Value = "Brenda Ribeiro"
$ws.Range("C8").Value = "Juridico"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45078
$ws.Range("G8").Value = 6528.28

$ws.Range("A9").Value = 43754
$ws.Range("B9").Value = "Sr. Gael Fogaça"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45085
$ws.Range("G9").Value = 5693.77

$ws.Range("A10").Value = 87852
$ws.Range("B10").Value = "Luna Cassiano"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45083
$ws.Range("G10").Value = 3906.77

$ws.Range("A11").Value = 90390
$ws.Range("B11").Value = "Juliana Azevedo"
$ws.Range("C11").Value = "Juridico"
$ws.Range("D11").Value = "Consulta medica"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45105
$ws.Range("G11").Value = 4563.78
